$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old title row (row 1), shifting all data up by one row
$ws.Rows.Item(1).Delete()

# Update header row text (capitalization / renaming) per the refreshed metadata
$headers = @(
    'Survey Year',
    'UnitID',
    'Institution Name',
    'Campus ID',
    'Campus Name',
    'Institution Size',
    'Murder/Non-Negligent Manslaughter',
    'Murder/Non-Negligent Manslaughter - Race',
    'Murder/Non-Negligent Manslaughter - Religion',
    'Murder/Non-Negligent Manslaughter - Sexual Orientation',
    'Murder/Non-Negligent Manslaughter - Gender',
    'Murder/Non-Negligent Manslaughter - Disability',
    'Murder/Non-Negligent Manslaughter - Ethnicity/National Origin',
    'Negligent Manslaughter',
    'Negligent Manslaughter - Race',
    'Negligent Manslaughter - Religion',
    'Negligent Manslaughter - Sexual Orientation',
    'Negligent Manslaughter - Gender',
    'Negligent Manslaughter - Disability',
    'Negligent Manslaughter - Ethnicity/National Origin',
    'Sex Offenses - Forcible',
    'Sex Offenses - Forcible - Race',
    'Sex Offenses - Forcible - Religion',
    'Sex Offenses - Forcible - Sexual Orientation',
    'Sex Offenses - Forcible - Gender',
    'Sex Offenses - Forcible - Disability',
    'Sex Offenses - Forcible - Ethnicity/National Origin',
    'Sex Offenses - Non-Forcible',
    'Sex Offenses - Non-Forcible -Race',
    'Sex Offenses - Non-Forcible - Religion',
    'Sex Offenses - Non-Forcible - Sexual Orientation',
    'Sex Offenses - Non-Forcible - Gender',
    'Sex Offenses - Non-Forcible - Disability',
    'Sex Offenses - Non-Forcible - Ethnicity/National Origin',
    'Robbery',
    'Robbery - Race',
    'Robbery - Religion',
    'Robbery - Sexual Orientation',
    'Robbery - Gender',
    'Robbery - Disability',
    'Robbery - Ethnicity/National Origin',
    'Aggravated Assault',
    'Aggravated Assault - Race',
    'Aggravated Assault - Religion',
    'Aggravated Assault - Sexual Orientation',
    'Aggravated Assault - Gender',
    'Aggravated Assault - Disability',
    'Aggravated Assault - Ethnicity/National Origin',
    'Burglary',
    'Burglary - Race',
    'Burglary - Religion',
    'Burglary - Sexual Orientation',
    'Burglary - Gender',
    'Burglary - Disability',
    'Burglary - Ethnicity/National Origin',
    'Motor Vehicle Theft',
    'Motor Vehicle Theft - Race',
    'Motor Vehicle Theft - Religion',
    'Motor Vehicle Theft - Sexual Orientation',
    'Motor Vehicle Theft - Gender',
    'Motor Vehicle Theft - Disability',
    'Motor Vehicle Theft - Ethnicity/National Origin',
    'Arson',
    'Arson - Race',
    'Arson - Religion',
    'Arson - Sexual Orientation',
    'Arson - Gender',
    'Arson - Disability',
    'Arson - Ethnicity/National Origin',
    'Simple Assault',
    'Simple Assault - Race',
    'Simple Assault - Religion',
    'Simple Assault - Sexual Orientation',
    'Simple Assault - Gender',
    'Simple Assault - Disability',
    'Simple Assault - Ethnicity/National Origin',
    'Larceny-Theft',
    'Larceny-Theft - Race',
    'Larceny-Theft - Religion',
    'Larceny-Theft - Sexual Orientation',
    'Larceny-Theft - Gender',
    'Larceny-Theft - Disability',
    'Larceny-Theft - Ethnicity/National Origin',
    'Intimidation',
    'Intimidation - Race',
    'Intimidation - Religion',
    'Intimidation - Sexual Orientation',
    'Intimidation - Gender',
    'Intimidation - Disability',
    'Intimidation - Ethnicity/National Origin',
    'Destruction/Damage/Vandalism of Property',
    'Destruction/Damage/Vandalism of Property - Race',
    'Destruction/Damage/Vandalism of Property - Religion',
    'Destruction/Damage/Vandalism of Property - Sexual Orientation',
    'Destruction/Damage/Vandalism of Property - Gender',
    'Destruction/Damage/Vandalism of Property - Disability',
    'Destruction/Damage/Vandalism of Property - Ethnicity/National Origin'
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
